$p = $ppt.ActivePresentation
$d = $p.Designs.Item(1)
try { Write-Host ("Design.Name: " + $d.Name) } catch { Write-Host ("ERR: " + $_.Exception.Message) }
try { $d.Name = "Integral2"; Write-Host ("set ok, now: " + $d.Name) } catch { Write-Host ("ERR set: " + $_.Exception.Message) }
